$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C26:I26").Value = 5

$ws.Range("H26").Select
